# Update gh-pages to output generated at 456a3b4
# Apply updated "想去人数" (want-to-go count) values across sheets.

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 6866
$ws1.Range("F8").Value = 6408
$ws1.Range("F10").Value = 1968
$ws1.Range("F17").Value = 52
$ws1.Range("F18").Value = 8099

# Sheet: 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 695

# Sheet: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 695
$ws4.Range("F9").Value = 6867
$ws4.Range("F12").Value = 6408
$ws4.Range("F14").Value = 1968
$ws4.Range("F22").Value = 52
$ws4.Range("F23").Value = 8099
